# Update backlog: remove unused "na" shared string, assign sprint numbers
# (as numeric values) to the "Assignée au sprint" column, scroll the view
# back to the top-left and move the selection to G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Sprint numbers for rows 2..18 (column G, "Assignée au sprint"),
# replacing the placeholder text "na" with real numeric sprint values.
$sprintValues = @(3, 4, 3, 3, 5, 5, 4, 4, 6, 4, 6, 5, 4, 6, 4, 6, 5)

for ($i = 0; $i -lt $sprintValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $sprintValues[$i]
}

# Reset the view: scroll back to the top and select G1.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("G1").Select()

$wb.Save()
